$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update submitNum (column J) and acceptNum (column K) values for the affected rows.
$ws.Range("J2").Value = 1358
$ws.Range("K2").Value = 737

$ws.Range("J3").Value = 1637

$ws.Range("J4").Value = 605
$ws.Range("K4").Value = 226

$ws.Range("J5").Value = 423
$ws.Range("K5").Value = 201

$ws.Range("J9").Value = 513

$ws.Range("J11").Value = 401
$ws.Range("K11").Value = 155

$ws.Range("J13").Value = 304
$ws.Range("K13").Value = 204

$ws.Range("J14").Value = 273

$ws.Range("J19").Value = 116
$ws.Range("K19").Value = 66

$ws.Range("J21").Value = 189
$ws.Range("K21").Value = 97

$ws.Range("J24").Value = 82

$ws.Range("J50").Value = 35
$ws.Range("K50").Value = 5

$ws.Range("J60").Value = 190
$ws.Range("K60").Value = 75
